$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the profile rows (2-5) so that:
#   new row 2 = Лингвистика (previously row 3)
#   new row 3 = Физика      (previously row 5)
#   new row 4 = Медицина    (previously row 2)
#   new row 5 = Математика  (previously row 4)
# Set all cell values explicitly to their final state.

$ws.Range("A2").Value = "Лингвистика"
$ws.Range("B2").Value = 0.0
$ws.Range("C2").Value = 0.0
$ws.Range("D2").Value = 1.0
$ws.Range("E2").Value = "Воронежский Литературно-Переводческий Университет; "

$ws.Range("A3").Value = "Физика"
$ws.Range("B3").Value = 4.5
$ws.Range("C3").Value = 8.0
$ws.Range("D3").Value = 2.0
$ws.Range("E3").Value = "Московский Выдуманный Университет; Московский Придуманный Институт; "

$ws.Range("A4").Value = "Медицина"
$ws.Range("B4").Value = 4.300000190734863
$ws.Range("C4").Value = 3.0
$ws.Range("D4").Value = 3.0
$ws.Range("E4").Value = "Московский Государственный Медицинский Университет; Тамбовский Университет Медицины; Самарский Медицинский Институт; "

$ws.Range("A5").Value = "Математика"
$ws.Range("B5").Value = 0.0
$ws.Range("C5").Value = 0.0
$ws.Range("D5").Value = 1.0
$ws.Range("E5").Value = "Казанский Университет Вычислений; "
